$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.169.69"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "3.086.79"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.86"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.09"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.37%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.079.41"
$ws.Range("E8").Value = "  -0.88%  "
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("E10").Value = "  +4.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.62"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("E12").Value = "  -2.87%  "
$ws.Range("E13").Value = "  -1.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.55"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.74%  "
$ws.Range("E15").Value = "  -1.27%  "
$ws.Range("D16").Value = "3.599.65"
$ws.Range("E16").Value = "  -0.84%  "
$ws.Range("D17").Value = "63.092.54"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("D19").Value = "3.083.89"
$ws.Range("E19").Value = "  -0.87%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "460.03"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.18"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.61%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.723"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.43"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.44%  "
$ws.Range("E24").Value = "  -3.14%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.00"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.43%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.12"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.65%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.88"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +7.30%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("E30").Value = "  -0.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.20"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.80"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.57"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.109"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.66%  "
$ws.Range("D35").Value = "0.0₃0840"
$ws.Range("E35").Value = "  -0.93%  "
$ws.Range("E36").Value = "  -1.18%  "
$ws.Range("E37").Value = "  -3.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.31"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.97"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.24"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.52%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "431.70"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.74"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("E43").Value = "  -0.61%  "
$ws.Range("D44").Value = "2.857.32"
$ws.Range("E44").Value = "  -2.53%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.268"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.33%  "
$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.108"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.08"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +2.40%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.92"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.110"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.99"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.04%  "
